# insert data into multiple worksheets in excel
$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename existing sheet and replace its data ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Patient Records"
$ws1.Cells.Clear()

$ws1.Range("A1").Value = "S.No."
$ws1.Range("B1").Value = "Patient No"
$ws1.Range("C1").Value = "First Name"
$ws1.Range("D1").Value = "Last Name"
$ws1.Range("E1").Value = "Condition"

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = "Dummy"
$ws1.Range("D2").Value = "Data"
$ws1.Range("E2").Value = "None"

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = 1
$ws1.Range("C3").Value = "test"
$ws1.Range("D3").Value = "test"
$ws1.Range("E3").Value = "heart failure"

# --- Sheet 2: new "Medications" worksheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Medications"

$ws2.Range("A1").Value = "S.No."
$ws2.Range("B1").Value = "Patient No."
$ws2.Range("C1").Value = "Medication"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = "blood-thinner"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = "thomapyrin"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = "Aspirin"

# --- Sheet 3: new "Activities" worksheet ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Activities"

$ws3.Range("A1").Value = "S.No."
$ws3.Range("B1").Value = "Patient No."
$ws3.Range("C1").Value = "Activities"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = 1
$ws3.Range("C2").Value = "blood-pressure"

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = "Pulse-rate"

$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = "blood-sugar-levels"

$ws3.Range("A5").Value = 4
$ws3.Range("B5").Value = 1
$ws3.Range("C5").Value = "Blood-drawn"

# restore per-sheet selections to match the final view state
$ws2.Range("C15").Select() | Out-Null
$ws3.Range("C8").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D1:D1048576").Select() | Out-Null
